$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 10393
$ws.Range("E2").Value = 521
$ws.Range("F2").Value = 521
$ws.Range("G2").Value = 346
$ws.Range("H2").Value = 160
$ws.Range("I2").Value = 167
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 12549
$ws.Range("L2").Value = 8637
$ws.Range("M2").Value = 3912
$ws.Range("N2").Value = 3819
$ws.Range("O2").Value = 93
$ws.Range("P2").Value = 106
$ws.Range("Q2").Value = 1395
$ws.Range("R2").Value = -1812
$ws.Range("S2").Value = 586
$ws.Range("T2").Value = 1458
$ws.Range("U2").Value = -63
$ws.Range("V2").Value = 4676
$ws.Range("W2").Value = 5.01
$ws.Range("X2").Value = 1.54
$ws.Range("Y2").Value = 4.42
$ws.Range("Z2").Value = 1.36
$ws.Range("AA2").Value = 220.8
$ws.Range("AB2").Value = 3539.92
$ws.Range("AC2").Value = 787
$ws.Range("AD2").Value = 68.59999999999999
$ws.Range("AE2").Value = 18048
$ws.Range("AF2").Value = 2.99
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 0.65
$ws.Range("AI2").Value = 44.46
$ws.Range("AJ2").Value = 21161313

# Row 3
$ws.Range("D3").Value = 11935
$ws.Range("E3").Value = 669
$ws.Range("F3").Value = 669
$ws.Range("G3").Value = 728
$ws.Range("H3").Value = 522
$ws.Range("I3").Value = 519
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 14176
$ws.Range("L3").Value = 9724
$ws.Range("M3").Value = 4452
$ws.Range("N3").Value = 4211
$ws.Range("O3").Value = 241
$ws.Range("P3").Value = 106
$ws.Range("Q3").Value = 1517
$ws.Range("R3").Value = -2273
$ws.Range("S3").Value = 1065
$ws.Range("T3").Value = 1919
$ws.Range("U3").Value = -403
$ws.Range("V3").Value = 5793
$ws.Range("W3").Value = 5.61
$ws.Range("X3").Value = 4.37
$ws.Range("Y3").Value = 12.94
$ws.Range("Z3").Value = 3.91
$ws.Range("AA3").Value = 218.43
$ws.Range("AB3").Value = 3948.38
$ws.Range("AC3").Value = 2455
$ws.Range("AD3").Value = 51.33
$ws.Range("AE3").Value = 19899
$ws.Range("AF3").Value = 6.33
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 0.28
$ws.Range("AI3").Value = 14.26
$ws.Range("AJ3").Value = 21161313

# Row 4
$ws.Range("D4").Value = 14322
$ws.Range("E4").Value = 703
$ws.Range("F4").Value = 703
$ws.Range("G4").Value = 180
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 123
$ws.Range("J4").Value = -67
$ws.Range("K4").Value = 25418
$ws.Range("L4").Value = 16213
$ws.Range("M4").Value = 9205
$ws.Range("N4").Value = 3736
$ws.Range("O4").Value = 5469
$ws.Range("P4").Value = 106
$ws.Range("Q4").Value = 1878
$ws.Range("R4").Value = -8354
$ws.Range("S4").Value = 7054
$ws.Range("T4").Value = 1913
$ws.Range("U4").Value = -35
$ws.Range("V4").Value = 10709
$ws.Range("W4").Value = 4.91
$ws.Range("X4").Value = 0.39
$ws.Range("Y4").Value = 3.1
$ws.Range("Z4").Value = 0.28
$ws.Range("AA4").Value = 176.13
$ws.Range("AB4").Value = 3967.39
$ws.Range("AC4").Value = 583
$ws.Range("AD4").Value = 120.74
$ws.Range("AE4").Value = 17654
$ws.Range("AF4").Value = 3.99
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 0.5
$ws.Range("AI4").Value = 60.03
$ws.Range("AJ4").Value = 21161313

# Row 5
$ws.Range("D5").Value = 17144
$ws.Range("E5").Value = 862
$ws.Range("F5").Value = 862
$ws.Range("G5").Value = 115
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = -14
$ws.Range("J5").Value = 114
$ws.Range("K5").Value = 24601
$ws.Range("L5").Value = 16824
$ws.Range("M5").Value = 7777
$ws.Range("N5").Value = 3080
$ws.Range("O5").Value = 4697
$ws.Range("P5").Value = 106
$ws.Range("Q5").Value = 1943
$ws.Range("R5").Value = -2859
$ws.Range("S5").Value = 360
$ws.Range("T5").Value = 2782
$ws.Range("U5").Value = -839
$ws.Range("V5").Value = 10617
$ws.Range("W5").Value = 5.03
$ws.Range("X5").Value = 0.59
$ws.Range("Y5").Value = -0.41
$ws.Range("Z5").Value = 0.4
$ws.Range("AA5").Value = 216.35
$ws.Range("AB5").Value = 3873.67
$ws.Range("AC5").Value = -66
$ws.Range("AD5").Value = -1119.52
$ws.Range("AE5").Value = 14554
$ws.Range("AF5").Value = 5.1
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 0.47
$ws.Range("AI5").Value = -528.08
$ws.Range("AJ5").Value = 21161313

# Row 6
$ws.Range("D6").Value = 17694
$ws.Range("E6").Value = 777
$ws.Range("F6").Value = 777
$ws.Range("G6").Value = -2105
$ws.Range("H6").Value = -1885
$ws.Range("I6").Value = -1407
$ws.Range("K6").Value = 22342
$ws.Range("L6").Value = 16839
$ws.Range("M6").Value = 5503
$ws.Range("N6").Value = 1118
$ws.Range("P6").Value = 106
$ws.Range("Q6").Value = 1087
$ws.Range("R6").Value = -434
$ws.Range("S6").Value = -58
$ws.Range("T6").Value = 2453
$ws.Range("U6").Value = -1366
$ws.Range("V6").Value = 9253
$ws.Range("W6").Value = 4.39
$ws.Range("X6").Value = -10.66
$ws.Range("Y6").Value = -67.01000000000001
$ws.Range("Z6").Value = -8.029999999999999
$ws.Range("AA6").Value = 306.01
$ws.Range("AB6").Value = 2566.82
$ws.Range("AC6").Value = -6647
$ws.Range("AD6").Value = -6.17
$ws.Range("AE6").Value = 5284
$ws.Range("AF6").Value = 7.76
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = -3.01
$ws.Range("AJ6").Value = 21161313

# Row 7
$ws.Range("D7").Value = 19503
$ws.Range("E7").Value = 1111
$ws.Range("G7").Value = -185
$ws.Range("H7").Value = -207
$ws.Range("I7").Value = -122
$ws.Range("K7").Value = 41768
$ws.Range("L7").Value = 36011
$ws.Range("M7").Value = 5757
$ws.Range("N7").Value = 1830
$ws.Range("P7").Value = 109
$ws.Range("Q7").Value = 4036
$ws.Range("R7").Value = -3322
$ws.Range("S7").Value = 915
$ws.Range("T7").Value = 3636
$ws.Range("U7").Value = -962
$ws.Range("W7").Value = 5.69
$ws.Range("X7").Value = -1.06
$ws.Range("Y7").Value = -8.24
$ws.Range("Z7").Value = -0.65
$ws.Range("AA7").Value = 625.54
$ws.Range("AC7").Value = -574
$ws.Range("AD7").Value = -45.72
$ws.Range("AE7").Value = 8647
$ws.Range("AF7").Value = 3.04
$ws.Range("AG7").Value = 213
$ws.Range("AH7").Value = 0.8100000000000001
$ws.Range("AI7").Value = -37.16

# Row 8
$ws.Range("D8").Value = 20951
$ws.Range("E8").Value = 1269
$ws.Range("G8").Value = 129
$ws.Range("H8").Value = 102
$ws.Range("I8").Value = 83
$ws.Range("K8").Value = 42487
$ws.Range("L8").Value = 36395
$ws.Range("M8").Value = 6093
$ws.Range("N8").Value = 2029
$ws.Range("P8").Value = 109
$ws.Range("Q8").Value = 4708
$ws.Range("R8").Value = -3030
$ws.Range("S8").Value = -114
$ws.Range("T8").Value = 2405
$ws.Range("U8").Value = 2209
$ws.Range("W8").Value = 6.06
$ws.Range("X8").Value = 0.49
$ws.Range("Y8").Value = 4.32
$ws.Range("Z8").Value = 0.24
$ws.Range("AA8").Value = 597.33
$ws.Range("AC8").Value = 394
$ws.Range("AD8").Value = 66.62
$ws.Range("AE8").Value = 9588
$ws.Range("AF8").Value = 2.74
$ws.Range("AG8").Value = 220
$ws.Range("AH8").Value = 0.84
$ws.Range("AI8").Value = 55.84

# Row 9
$ws.Range("D9").Value = 22498
$ws.Range("E9").Value = 1442
$ws.Range("G9").Value = 362
$ws.Range("H9").Value = 297
$ws.Range("I9").Value = 219
$ws.Range("K9").Value = 44170
$ws.Range("L9").Value = 37662
$ws.Range("M9").Value = 6508
$ws.Range("N9").Value = 2207
$ws.Range("P9").Value = 109
$ws.Range("Q9").Value = 4743
$ws.Range("R9").Value = -3106
$ws.Range("S9").Value = 187
$ws.Range("T9").Value = 2335
$ws.Range("U9").Value = 2146
$ws.Range("W9").Value = 6.41
$ws.Range("X9").Value = 1.32
$ws.Range("Y9").Value = 10.36
$ws.Range("Z9").Value = 0.6899999999999999
$ws.Range("AA9").Value = 578.6900000000001
$ws.Range("AC9").Value = 1037
$ws.Range("AD9").Value = 25.32
$ws.Range("AE9").Value = 10432
$ws.Range("AF9").Value = 2.52
$ws.Range("AG9").Value = 231
$ws.Range("AH9").Value = 0.88
$ws.Range("AI9").Value = 22.25
